$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 16 (Bradley Beal) and row 19 (Ja Morant)
$row16 = @($ws.Range("A16").Value(), $ws.Range("B16").Value(), $ws.Range("C16").Value())
$row19 = @($ws.Range("A19").Value(), $ws.Range("B19").Value(), $ws.Range("C19").Value())

$ws.Range("A16").Value = $row19[0]
$ws.Range("B16").Value = $row19[1]
$ws.Range("C16").Value = $row19[2]

$ws.Range("A19").Value = $row16[0]
$ws.Range("B19").Value = $row16[1]
$ws.Range("C19").Value = $row16[2]
